$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows for account 004493324 (DANIEL) and 004363250 (HELIO)
$ws.Rows("4:5").Delete()

# Update the balance for account 004472386 (GABRIEL), now shifted up to row 4
$ws.Range("C4").Value = 35000
